$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Farmacias")

# Append the new pharmacy row (row 17) to the Farmacias table
$ws.Range("A17").Value = "Farmacia Alem"
$ws.Range("B17").Value = "Leandro Alem 2654"
$ws.Range("C17").Value = "Ushuaia"
$ws.Range("D17").Value = "Ushuaia"
$ws.Range("E17").Value = "Tierra del Fuego"
$ws.Range("F17").Value = 2901584008
$ws.Range("G17").Value = "(2901)584-008"
$ws.Range("H17").Value = 2901425045
$ws.Range("I17").Value = "(2901) 425-045"

# Match the style used by the new row: left aligned, vertically centered
$ws.Range("F17").HorizontalAlignment = -4131
$ws.Range("F17").VerticalAlignment = -4108

# Match the style used for the existing phone columns (left aligned, like H3:H16)
$ws.Range("H17").HorizontalAlignment = -4131
$ws.Range("I17").HorizontalAlignment = -4131

# Grow the table (TablaFarmacias) so it includes the newly added row
$lo = $ws.ListObjects.Item("TablaFarmacias")
$lo.Resize($ws.Range("A1:J17"))

# Update selection left on sheet1 when the file was saved
$ws.Range("G26").Select()

# Leave "Instrucciones" as the active sheet when the file was saved
$ws2 = $wb.Worksheets.Item("Instrucciones")
$ws2.Activate()
